$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Bruno Cortes - Marketing Digital"
$ws.Range("B2").Value = "Estr. da Cachamorra, 1233 - Campo Grande, Rio de Janeiro - RJ, 23040-150"
$ws.Range("C2").Value = "(21) 98836-0205"
$ws.Range("D2").Value = "brunocortes.com.br"

$ws.Range("A3").Value = "Winner Digital Pro | Agencia de Marketing Digital | Gestor de Trafego | Criação de Sites | Social Media"
$ws.Range("B3").Value = "R. Gramado, 475 - Campo Grande, Rio de Janeiro - RJ, 23050-090"
$ws.Range("C3").Value = "(21) 98552-3425"
$ws.Range("D3").Value = "winnerdigitalpro.com"

$ws.Range("A4").Value = "RAES Marketing Digital"
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""

$ws.Range("A5").Value = "Bruno Cortes - Marketing Digital"
$ws.Range("B5").Value = "Estr. da Cachamorra, 1233 - Campo Grande, Rio de Janeiro - RJ, 23040-150"
$ws.Range("C5").Value = "(21) 98836-0205"
$ws.Range("D5").Value = "brunocortes.com.br"

$ws.Range("A6").Value = "E-mind Negócios Digitais"
$ws.Range("B6").Value = "Estr. da Cachamorra, 350 - Bloco 1 Sala 415 - Campo Grande, Rio de Janeiro - RJ, 23040-150"
$ws.Range("C6").Value = "(21) 96640-2049"
$ws.Range("D6").Value = "emindnegocios.com.br"

$ws.Range("A7").Value = "New Digital | Agencia de Marketing Digital | Gestor de Trafego | Criação de Sites"
$ws.Range("B7").Value = "Centro Comercial Business Completo, Av. Maria Teresa, 75 - Campo Grande, Rio de Janeiro - RJ, 23050-160"
$ws.Range("C7").Value = "(21) 99880-4831"
$ws.Range("D7").Value = "newdigitalpro.com.br"

$ws.Range("A8").Value = "Agência e Produtora Páginas e Aplicativos"
$ws.Range("B8").Value = "Office Mall - RJ - Estr. da Cachamorra, 350 - Bloco 3 - Sala 461 - Campo Grande, Rio de Janeiro - RJ, 23040-150"
$ws.Range("C8").Value = "(21) 98046-3733"
$ws.Range("D8").Value = ""

$ws.Range("A9").Value = "Agência de Marketing Digital - BeHype Media"
$ws.Range("B9").Value = "Av. Maria Teresa, 260 - Campo Grande, Rio de Janeiro - RJ, 23050-160"
$ws.Range("C9").Value = "(21) 97585-6922"
$ws.Range("D9").Value = "behypemedia.com"

$ws.Range("A10").Value = "Agência de Marketing Digital | Ouro Space"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "(21) 97274-1361"
$ws.Range("D10").Value = ""

$ws.Range("A11").Value = "TM Group - Agência de Marketing"
$ws.Range("B11").Value = "R. Vicente Perrota - Campo Grande, Rio de Janeiro - RJ, 23036-180"
$ws.Range("C11").Value = "(21) 97949-3723"
$ws.Range("D11").Value = "agenciatmgroup.com.br"
